# Generate Report for Handoff
#
# Adds two new "Ready for handoff" file entries to the localization-status
# report:
#   - 13ab9a5f-c888-40e6-96e7-24f92956d476.md
#   - 57c49093-ee5e-4a54-85ab-e2f1cfd93210.md
#
# They are inserted, on every sheet (Overview, zh-cn, de-de), right before
# the pre-existing "ac038f04-339a-42d5-a055-cd9711f4b8af.md" row and before
# the ".localization-config" row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ===========================================================================
# Overview sheet (columns A:C are always fully populated on every row, so a
# plain double row-insert at row 3 is safe - no "gap" columns to worry
# about).
# ===========================================================================
$ws1.Rows.Item(3).Insert()
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "13ab9a5f-c888-40e6-96e7-24f92956d476.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = "57c49093-ee5e-4a54-85ab-e2f1cfd93210.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

# (rows 5 = ac038f04..., 6 = .localization-config already hold the right
#  values - Insert() shifted them down intact.)

# ===========================================================================
# zh-cn / de-de sheets: columns E:F are only populated on row 2 (the
# "1b18accd" row). Inserting new rows directly at row 3 would copy row 2's
# formatting down and leave stray empty-but-styled E/F cells behind that
# can never be fully cleared again. Instead, insert at row 4 - which copies
# the clean A/B/C/D/G/H-only formatting of row 3 (the "ac038f04" row) - then
# relocate the "ac038f04" row's own data down into the freed row 5, and
# finally overwrite rows 3 and 4 with the two new entries.
# ===========================================================================
function Add-HandoffRows($ws, $zhOrDe, $dt) {
    $ws.Rows.Item(4).Insert()
    $ws.Rows.Item(4).Insert()

    # Relocate the existing "ac038f04" row (still sitting at row 3) down to
    # row 5, cell by cell so only the columns that really hold data get
    # touched (keeps the sparse E/F gap intact).
    foreach ($col in @("A","B","C","D","E","F","G","H","I")) {
        $val = $ws.Range($col + "3").Text
        if ($val -ne "") {
            $ws.Range($col + "5").Value = $val
        }
    }

    $ws.Range("A3").Value = "13ab9a5f-c888-40e6-96e7-24f92956d476.md"
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("C3").Value = "13ab9a5f-c888-40e6-96e7-24f92956d476.c408e3d2184f9bc536a94a0263f0ba59e24717b9." + $zhOrDe + ".xlf"
    $ws.Range("D3").Value = $dt
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Include"

    $ws.Range("A4").Value = "57c49093-ee5e-4a54-85ab-e2f1cfd93210.md"
    $ws.Range("B4").Value = "Ready for handoff"
    $ws.Range("C4").Value = "57c49093-ee5e-4a54-85ab-e2f1cfd93210.5c6b5fd06a52bd2a636cd029eb216b690fea975a." + $zhOrDe + ".xlf"
    $ws.Range("D4").Value = $dt
    $ws.Range("G4").Value = "0001-01-01 00:00:00"
    $ws.Range("H4").Value = "Include"
}

Add-HandoffRows $ws2 "zh-cn" "2016-03-09 00:49:45"
Add-HandoffRows $ws3 "de-de" "2016-03-09 00:49:54"

# ===========================================================================
# Hyperlinks do not automatically follow the rows they were anchored to when
# rows are inserted, so rebuild every hyperlink on every sheet from scratch
# once all cell values sit in their final positions.
# ===========================================================================
$ws1.Hyperlinks.Delete()
$ws2.Hyperlinks.Delete()
$ws3.Hyperlinks.Delete()

# -- Overview --
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c9c24e6afd7052ae63682abd53b82e9faa8b87e/e2e/1b18accd-ef67-4ed0-b431-a21bf8f620ba.md", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c408e3d2184f9bc536a94a0263f0ba59e24717b9/e2e/13ab9a5f-c888-40e6-96e7-24f92956d476.md", $null, $null, "13ab9a5f-c888-40e6-96e7-24f92956d476.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c6b5fd06a52bd2a636cd029eb216b690fea975a/e2e/57c49093-ee5e-4a54-85ab-e2f1cfd93210.md", $null, $null, "57c49093-ee5e-4a54-85ab-e2f1cfd93210.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6033ff7ce15ee4d1c2be51c6a10f43ef3a1d0aae/e2e/ac038f04-339a-42d5-a055-cd9711f4b8af.md", $null, $null, "ac038f04-339a-42d5-a055-cd9711f4b8af.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/7c9c24e6afd7052ae63682abd53b82e9faa8b87e/.localization-config", $null, $null, ".localization-config") | Out-Null

# -- zh-cn --
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c9c24e6afd7052ae63682abd53b82e9faa8b87e/e2e/1b18accd-ef67-4ed0-b431-a21bf8f620ba.md", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1cc8a437185c1077f4b19c83e9509c8a57f4a8ce/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.zh-cn.xlf", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/cf8567802f47071bc35a6f8b31235206b51c42ae/e2e/1b18accd-ef67-4ed0-b431-a21bf8f620ba.md", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2d4263819f8cee9fb7ce6924c3edc9fff12f81ae/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.zh-cn.xlf", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c408e3d2184f9bc536a94a0263f0ba59e24717b9/e2e/13ab9a5f-c888-40e6-96e7-24f92956d476.md", $null, $null, "13ab9a5f-c888-40e6-96e7-24f92956d476.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c408e3d2184f9bc536a94a0263f0ba59e24717b9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13ab9a5f-c888-40e6-96e7-24f92956d476.c408e3d2184f9bc536a94a0263f0ba59e24717b9.zh-cn.xlf", $null, $null, "13ab9a5f-c888-40e6-96e7-24f92956d476.c408e3d2184f9bc536a94a0263f0ba59e24717b9.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c6b5fd06a52bd2a636cd029eb216b690fea975a/e2e/57c49093-ee5e-4a54-85ab-e2f1cfd93210.md", $null, $null, "57c49093-ee5e-4a54-85ab-e2f1cfd93210.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c6b5fd06a52bd2a636cd029eb216b690fea975a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/57c49093-ee5e-4a54-85ab-e2f1cfd93210.5c6b5fd06a52bd2a636cd029eb216b690fea975a.zh-cn.xlf", $null, $null, "57c49093-ee5e-4a54-85ab-e2f1cfd93210.5c6b5fd06a52bd2a636cd029eb216b690fea975a.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6033ff7ce15ee4d1c2be51c6a10f43ef3a1d0aae/e2e/ac038f04-339a-42d5-a055-cd9711f4b8af.md", $null, $null, "ac038f04-339a-42d5-a055-cd9711f4b8af.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a03571732da0c14e439913c4bcb237257d4d2969/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ac038f04-339a-42d5-a055-cd9711f4b8af.24729ead3b959637028f29622d6ffbda2f5e36bc.zh-cn.xlf", $null, $null, "ac038f04-339a-42d5-a055-cd9711f4b8af.24729ead3b959637028f29622d6ffbda2f5e36bc.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/7c9c24e6afd7052ae63682abd53b82e9faa8b87e/.localization-config", $null, $null, ".localization-config") | Out-Null

# -- de-de --
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c9c24e6afd7052ae63682abd53b82e9faa8b87e/e2e/1b18accd-ef67-4ed0-b431-a21bf8f620ba.md", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43b451f06a5b45d69466f12fcab835cc9190598f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.de-de.xlf", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/40e88a0ebc6e6278ec658258c0581ba7540caebe/e2e/1b18accd-ef67-4ed0-b431-a21bf8f620ba.md", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8b780d76b5cdc8fc356a82dbdea29e5cbd582d91/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.de-de.xlf", $null, $null, "1b18accd-ef67-4ed0-b431-a21bf8f620ba.dc295c9205e5ec0c4998d1e4df546e981550c0eb.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c408e3d2184f9bc536a94a0263f0ba59e24717b9/e2e/13ab9a5f-c888-40e6-96e7-24f92956d476.md", $null, $null, "13ab9a5f-c888-40e6-96e7-24f92956d476.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c408e3d2184f9bc536a94a0263f0ba59e24717b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13ab9a5f-c888-40e6-96e7-24f92956d476.c408e3d2184f9bc536a94a0263f0ba59e24717b9.de-de.xlf", $null, $null, "13ab9a5f-c888-40e6-96e7-24f92956d476.c408e3d2184f9bc536a94a0263f0ba59e24717b9.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5c6b5fd06a52bd2a636cd029eb216b690fea975a/e2e/57c49093-ee5e-4a54-85ab-e2f1cfd93210.md", $null, $null, "57c49093-ee5e-4a54-85ab-e2f1cfd93210.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c6b5fd06a52bd2a636cd029eb216b690fea975a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/57c49093-ee5e-4a54-85ab-e2f1cfd93210.5c6b5fd06a52bd2a636cd029eb216b690fea975a.de-de.xlf", $null, $null, "57c49093-ee5e-4a54-85ab-e2f1cfd93210.5c6b5fd06a52bd2a636cd029eb216b690fea975a.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6033ff7ce15ee4d1c2be51c6a10f43ef3a1d0aae/e2e/ac038f04-339a-42d5-a055-cd9711f4b8af.md", $null, $null, "ac038f04-339a-42d5-a055-cd9711f4b8af.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43b451f06a5b45d69466f12fcab835cc9190598f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ac038f04-339a-42d5-a055-cd9711f4b8af.24729ead3b959637028f29622d6ffbda2f5e36bc.de-de.xlf", $null, $null, "ac038f04-339a-42d5-a055-cd9711f4b8af.24729ead3b959637028f29622d6ffbda2f5e36bc.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/7c9c24e6afd7052ae63682abd53b82e9faa8b87e/.localization-config", $null, $null, ".localization-config") | Out-Null

Write-Host "Done."
